$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Data: was Madagascar master data (fra / MDG); now Sierra Leone
# (eng / SLE). Also lower-case the resident-client service account.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = "SLE"
$ws.Range("C2").Value = "globaladmin"
$ws.Range("D2").Value = $true

$ws.Range("A3").Value = "eng"
$ws.Range("B3").Value = "SLE"
$ws.Range("C3").Value = "service-account-mosip-resident-client"
$ws.Range("D3").Value = $true

# ------------------------------------------------------------------
# Header row (row 1): drop wrap-text, thin out the border, make the
# "zone_code" header match the bold "lang_code" header style, and
# mark the last two header cells as Text-formatted.
# ------------------------------------------------------------------
$ws.Range("A1:D1").WrapText = $false
$ws.Range("A1:D1").Borders.Weight = 2
$ws.Range("A1:D1").Borders.ColorIndex = 1

$hdr2 = $ws.Range("B1")
$hdr2.Font.Name = "Cambria"
$hdr2.Font.Bold = $true
$hdr2.HorizontalAlignment = -4108
$hdr2.VerticalAlignment = -4160

$ws.Range("C1:D1").NumberFormat = "@"

# ------------------------------------------------------------------
# Data rows (2-3): plain (unstyled) text cells for lang/zone/usr_id,
# italic smaller font for the zone_code column, and a custom
# TRUE/FALSE display format for the boolean column.
# ------------------------------------------------------------------
$rA = $ws.Range("A2:A3")
$rA.Font.Name = "Calibri"
$rA.Font.Size = 11
$rA.Font.Italic = $false
$rA.Borders.LineStyle = 0
$rA.WrapText = $false

$rC = $ws.Range("C2:C3")
$rC.Font.Name = "Calibri"
$rC.Font.Size = 11
$rC.Font.Italic = $false
$rC.Borders.LineStyle = 0
$rC.WrapText = $false

$rB = $ws.Range("B2:B3")
$rB.Style = "Normal"
$rB.Font.Italic = $true
$rB.Font.Size = 10

$rD = $ws.Range("D2:D3")
$rD.Style = "Normal"
$rD.NumberFormat = '"TRUE";"TRUE";"FALSE"'

# ------------------------------------------------------------------
# Rows: drop the explicit (thick-bottom) row heights back to default.
# ------------------------------------------------------------------
$ws.Rows("1:3").AutoFit()

# ------------------------------------------------------------------
# Columns: narrower A/B, wider C, new D column width.
# ------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 9.5
$ws.Columns(2).ColumnWidth = 10
$ws.Columns(3).ColumnWidth = 31.5
$ws.Columns(4).ColumnWidth = 8.166666666666666

# ------------------------------------------------------------------
# View: selection cursor moved to C8.
# ------------------------------------------------------------------
$ws.Range("C8").Select()
